$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new week's driving data (Lyft) as row 17
$ws.Range("A17").Value = "Lyft"
$ws.Range("B17").Value = 44360
$ws.Range("C17").Value = 57

$ws.Range("D17").Formula = "=27*60"
$ws.Range("E17").Formula = "=30*60"

$ws.Range("F17").Value = 536.59
$ws.Range("G17").Value = 47.64
$ws.Range("H17").Value = 72.78
$ws.Range("I17").Value = 17.77
$ws.Range("J17").Value = 19.58
$ws.Range("K17").Value = 255
$ws.Range("L17").Value = 0

$ws.Range("M17").Formula = "=SUM(F17:L17)"
$ws.Range("N17").Formula = "=M17-J17"

# Copy the formatting from the row above (dates, currency, centered alignment)
# so the new row matches the existing table styling exactly.
$ws.Range("A16:N16").Copy()
$ws.Range("A17:N17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the active selection to where the user clicked next (below the new row)
$ws.Range("F18").Select()
